$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 434; existing rows 434..487 shift down to 435..488,
# and the sheet dimension grows from A1:R487 to A1:R488.
$ws.Rows(434).Insert()

# Populate the newly inserted row 434 with the new weekly data point
# (same market/category/variety metadata as the surrounding rows).
$ws.Cells.Item(434, 1).Value  = 4
$ws.Cells.Item(434, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(434, 3).Value  = "Los Lagos"
$ws.Cells.Item(434, 4).Value  = 45127
$ws.Cells.Item(434, 5).Value  = 10
$ws.Cells.Item(434, 6).Value  = 100112037
$ws.Cells.Item(434, 7).Value  = "Cebollín"
$ws.Cells.Item(434, 8).Value  = "Sin especificar"
$ws.Cells.Item(434, 9).Value  = "Primera"
$ws.Cells.Item(434, 10).Value = 70
$ws.Cells.Item(434, 11).Value = 6000
$ws.Cells.Item(434, 12).Value = 6500
$ws.Cells.Item(434, 13).Value = 6250
$ws.Cells.Item(434, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(434, 15).Value = "Región Metropolitana"
$ws.Cells.Item(434, 16).Value = 174
$ws.Cells.Item(434, 17).Value = 36
$ws.Cells.Item(434, 18).Value = "Hortaliza"
